# DataValidationRules.xlsx - "Show UI when List and Regex rules fail"
#
# Replaces the old Doors/Floors validation-rule sample rows with a new
# List + Regex example pair that demonstrates user-visible failure
# messages, and trims the trailing padding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Wipe the old sample data (rows 2-6, cols A-F) - contents AND formats
#    so stale styles/shared-strings don't linger.
# ---------------------------------------------------------------------
$ws.Range("A2:F6").Clear()

# ---------------------------------------------------------------------
# 2. Row 2 - "Floors / Comments / List" example
#    (cells are written left-to-right so first-use order in the saved
#    sharedStrings table lines up with the authored workbook)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Floors"
$ws.Range("B2").Value = "Comments"
$ws.Range("C2").Value = "List"
$ws.Range("D2").Value = "the,quick,brown,fox"
$ws.Range("E2").Value = "Y"
$ws.Range("F2").Value = "Enter a value from " + [string][char]0x201C + "the quick brown fox" + [string][char]0x201D

$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 3. Row 3 - "<all> / Mark / Regex" example
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "<all>"
$ws.Range("B3").Value = "Mark"
$ws.Range("C3").Value = "Regex"
$ws.Range("D3").Value = "^\d+$"
$ws.Range("E3").Value = "N"
$ws.Range("F3").Value = "Mark must contain only numbers"

$ws.Range("D3").WrapText = $true

# ---------------------------------------------------------------------
# 4. Rows 4 & 5 stay blank (already cleared above).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5. Row 6 - two empty, word-wrapped placeholder cells.
# ---------------------------------------------------------------------
$ws.Range("D6").WrapText = $true
$ws.Range("E6").WrapText = $true

# ---------------------------------------------------------------------
# 6. Drop the four now-unused padding rows at the bottom (1003-1006),
#    shifting everything below up so the sheet ends at row 1002.
# ---------------------------------------------------------------------
$ws.Rows("1003:1006").Delete()

# ---------------------------------------------------------------------
# 7. Move the active selection to E5, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("E5").Select() | Out-Null
